# The document originally splits several headline strings into many
# single-word runs (one <w:r> per word/space). This collapses each of
# those runs back into one contiguous run per paragraph by doing a
# Find/Replace across the whole (fragmented) text with the same text as
# the replacement - Word rebuilds the matched range as a single run.

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2

# Title paragraph: "Questions: Introduction to complex numbers"
$d.Content.Find.Execute("Questions: Introduction to complex numbers", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Questions: Introduction to complex numbers", 2)

# Author paragraph: "Tom Coleman"
$d.Content.Find.Execute("Tom Coleman", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Tom Coleman", 2)

# Abstract paragraph
$d.Content.Find.Execute("A selection of questions for the study guide on introduction to complex numbers.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "A selection of questions for the study guide on introduction to complex numbers.", 2)
